$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing rows 8, 11, 15, 16 need to move to 9, 17, 21, 22.
# Clear their old contents first (since they are moving to new rows),
# then rewrite all cell values/formulas in top-to-bottom, left-to-right
# order so that new shared strings get appended in the right order.

$ws.Range("A8").ClearContents()
$ws.Range("A11").ClearContents()
$ws.Range("A15").ClearContents()
$ws.Range("A16").ClearContents()

# Row 1
$ws.Range("A1").Value = "Budget for dinner event"

# Row 4
$ws.Range("A4").Value = "cookbook production cost, maybe should be broken down into some parts"

# Row 5
$ws.Range("A5").Value = "150 copies after we have the art, by CCSU Copy Center"
$ws.Range("F5").Value = 1000

# Row 6
$ws.Range("A6").Value = "the art: photography costs"

# Row 7
$ws.Range("A7").Value = "the recipes: "

# Row 9
$ws.Range("A9").Value = "venue use, insurance, security?"

# Row 10
$ws.Range("A10").Value = "4 hours of the site"
$ws.Range("F10").Value = 500

# Row 11
$ws.Range("A11").Value = "insurance: self insured"

# Row 12
$ws.Range("A12").Value = "security:"
$ws.Range("B12").Value = "none"

# Row 17
$ws.Range("A17").Value = "per person dining cost"

# Row 18
$ws.Range("A18").Value = 150
$ws.Range("B18").Value = "people"
$ws.Range("C18").Value = 25
$ws.Range("D18").Value = "per seat"
$ws.Range("F18").Formula = "=A18*C18"

# Row 21
$ws.Range("A21").Value = "calculate how many people we may serve"

# Row 22
$ws.Range("A22").Value = "make sure the room seats this many"

# Update the active selection cell to match the target state
$ws.Range("A24").Select()
